$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the value to be stored as text even if it looks numeric,
    # preserving the original (default) cell style.
    $style = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $style
}

$ws.Range('D2').Value = '84.709.16'
$ws.Range('E2').Value = '  +5.73%  '
$ws.Range('D3').Value = '3.304.54'
$ws.Range('E3').Value = '  +2.48%  '
Set-TextValue $ws.Range('D4') '0.999'
$ws.Range('E4').Value = '  +0.07%  '
Set-TextValue $ws.Range('D5') '219.25'
$ws.Range('E5').Value = '  +3.29%  '
Set-TextValue $ws.Range('D6') '636.39'
$ws.Range('E6').Value = '  -0.80%  '
Set-TextValue $ws.Range('D7') '0.325'
$ws.Range('E7').Value = '  +20.04%  '
Set-TextValue $ws.Range('D8') '0.998'
$ws.Range('E8').Value = '  -0.01%  '
Set-TextValue $ws.Range('D9') '0.593'
$ws.Range('E9').Value = '  -2.27%  '
$ws.Range('D10').Value = '3.303.30'
$ws.Range('E10').Value = '  +2.55%  '
Set-TextValue $ws.Range('D11') '0.596'
$ws.Range('E11').Value = '  -4.09%  '
Set-TextValue $ws.Range('D12') '0.0000279'
$ws.Range('E12').Value = '  +1.80%  '
$ws.Range('E13').Value = '  +0.05%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue $ws.Range('D14') '34.27'
$ws.Range('E14').Value = '  +4.35%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '3.899.43'
$ws.Range('E15').Value = '  +2.40%  '
Set-TextValue $ws.Range('D16') '5.44'
$ws.Range('E16').Value = '  -0.20%  '
$ws.Range('D17').Value = '84.522.99'
$ws.Range('E17').Value = '  +5.87%  '
$ws.Range('D18').Value = '3.284.95'
$ws.Range('E18').Value = '  +2.25%  '
$ws.Range('B19').Value = 'SuiNetwork'
$ws.Range('C19').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue $ws.Range('D19') '3.23'
$ws.Range('E19').Value = '  +6.27%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range('D20') '14.62'
$ws.Range('E20').Value = '  -0.60%  '
Set-TextValue $ws.Range('D21') '9.22'
$ws.Range('E21').Value = '  -2.22%  '
Set-TextValue $ws.Range('D22') '436.36'
$ws.Range('E22').Value = '  -2.80%  '
Set-TextValue $ws.Range('D23') '5.23'
$ws.Range('E23').Value = '  -2.48%  '
Set-TextValue $ws.Range('D24') '7.46'
$ws.Range('E24').Value = '  +4.02%  '
Set-TextValue $ws.Range('D25') '5.55'
$ws.Range('E25').Value = '  +13.37%  '
Set-TextValue $ws.Range('D26') '12.22'
$ws.Range('E26').Value = '  +11.21%  '
$ws.Range('D27').Value = '3.446.00'
$ws.Range('E27').Value = '  +2.37%  '
Set-TextValue $ws.Range('D28') '78.08'
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('E29').Value = '  +3.08%  '
Set-TextValue $ws.Range('D30') '1.00'
$ws.Range('E30').Value = '  -0.20%  '
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range('D31') '599.17'
$ws.Range('E31').Value = '  +5.40%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range('D32') '9.29'
$ws.Range('E32').Value = '  -0.45%  '
$ws.Range('B33').Value = 'Cronos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range('D33') '0.162'
$ws.Range('E33').Value = '  +31.14%  '
$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range('D34') '1.58'
$ws.Range('E34').Value = '  +2.65%  '
$ws.Range('B35').Value = 'Binance-PegBSC-USD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue $ws.Range('D35') '1.00'
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('E36').Value = '  -2.41%  '
Set-TextValue $ws.Range('D37') '2.04'
$ws.Range('E37').Value = '  -0.92%  '
Set-TextValue $ws.Range('D38') '23.28'
$ws.Range('E38').Value = '  -1.14%  '
Set-TextValue $ws.Range('D39') '6.43'
$ws.Range('E39').Value = '  +9.83%  '
$ws.Range('E40').Value = '  -0.01%  '
Set-TextValue $ws.Range('D41') '0.416'
$ws.Range('E41').Value = '  -0.13%  '
Set-TextValue $ws.Range('D42') '3.14'
$ws.Range('E42').Value = '  +13.59%  '
Set-TextValue $ws.Range('D43') '2.06'
$ws.Range('E43').Value = '  +11.71%  '
Set-TextValue $ws.Range('D44') '20.96'
$ws.Range('E44').Value = '  +3.17%  '
Set-TextValue $ws.Range('D45') '158.52'
$ws.Range('E45').Value = '  -3.22%  '
Set-TextValue $ws.Range('D47') '190.67'
$ws.Range('E47').Value = '  -1.33%  '
Set-TextValue $ws.Range('D48') '45.22'
$ws.Range('E48').Value = '  +4.74%  '
$ws.Range('E49').Value = '  -0.08%  '
Set-TextValue $ws.Range('D50') '0.789'
$ws.Range('E50').Value = '  -1.87%  '
Set-TextValue $ws.Range('D51') '26.76'
$ws.Range('E51').Value = '  +2.77%  '
